$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 557771.1
$ws.Range("I40").Value = 1599.875
$ws.Range("J40").Value = 1002708.1
$ws.Range("K40").Value = 1599.875
$ws.Range("L40").Value = 1002708.1
$ws.Range("M40").Value = -1424.875
$ws.Range("N40").Value = -1003058.1
$ws.Range("H80").Value = 357.6111
$ws.Range("I80").Value = 329
$ws.Range("J80").Value = 500.66666
$ws.Range("K80").Value = 987
$ws.Range("L80").Value = 1501.99998
$ws.Range("M80").Value = 11
$ws.Range("N80").Value = -3497.99998
$ws.Range("H83").Value = 357.6111
$ws.Range("I83").Value = 329
$ws.Range("J83").Value = 500.66666
$ws.Range("K83").Value = 2961
$ws.Range("L83").Value = 4505.99994
$ws.Range("M83").Value = 2031
$ws.Range("N83").Value = -14489.99994
$ws.Range("H88").Value = 30429780
$ws.Range("J88").Value = 36515136
$ws.Range("L88").Value = 36515136
$ws.Range("N88").Value = -36515948
$ws.Range("H91").Value = 30429780
$ws.Range("J91").Value = 36515136
$ws.Range("L91").Value = 36515136
$ws.Range("N91").Value = -36517944
$ws.Range("H135").Value = 951.77356
$ws.Range("I135").Value = 756.5333000000001
$ws.Range("K135").Value = 6808.7997
$ws.Range("M135").Value = -4273.7997

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1098.8
$ws.Range("I61").Value = 959.02856
$ws.Range("J61").Value = 1424.9333
$ws.Range("K61").Value = 959.02856
$ws.Range("L61").Value = 1424.9333
$ws.Range("M61").Value = -747.02856
$ws.Range("N61").Value = -1848.9333
$ws.Range("H74").Value = 723.7778
$ws.Range("I74").Value = 727.7895
$ws.Range("K74").Value = 727.7895
$ws.Range("M74").Value = 146.2105
$ws.Range("H77").Value = 723.7778
$ws.Range("I77").Value = 727.7895
$ws.Range("K77").Value = 3638.9475
$ws.Range("M77").Value = 729.0525000000002
$ws.Range("H97").Value = 878.65717
$ws.Range("I97").Value = 741.2692
$ws.Range("J97").Value = 1275.5555
$ws.Range("K97").Value = 741.2692
$ws.Range("L97").Value = 1275.5555
$ws.Range("M97").Value = -245.2692
$ws.Range("N97").Value = -2267.5555
$ws.Range("H122").Value = 1244.5714
$ws.Range("I122").Value = 978
$ws.Range("K122").Value = 2934
$ws.Range("M122").Value = -484
$ws.Range("H136").Value = 1098.8
$ws.Range("I136").Value = 959.02856
$ws.Range("J136").Value = 1424.9333
$ws.Range("K136").Value = 2877.08568
$ws.Range("L136").Value = 4274.7999
$ws.Range("M136").Value = -327.0856800000001
$ws.Range("N136").Value = -9374.7999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5056931.5
$ws.Range("I134").Value = 2012.75
$ws.Range("J134").Value = 7945456
$ws.Range("K134").Value = 6038.25
$ws.Range("L134").Value = 23836368
$ws.Range("M134").Value = -3503.25
$ws.Range("N134").Value = -23841438

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1284.6938
$ws.Range("I31").Value = 960.13043
$ws.Range("J31").Value = 1571.8077
$ws.Range("K31").Value = 960.13043
$ws.Range("L31").Value = 1571.8077
$ws.Range("M31").Value = -665.13043
$ws.Range("N31").Value = -2161.8077
$ws.Range("H34").Value = 1284.6938
$ws.Range("I34").Value = 960.13043
$ws.Range("J34").Value = 1571.8077
$ws.Range("K34").Value = 960.13043
$ws.Range("L34").Value = 1571.8077
$ws.Range("M34").Value = -758.13043
$ws.Range("N34").Value = -1975.8077
$ws.Range("H99").Value = 27032594
$ws.Range("I99").Value = 41673200
$ws.Range("J99").Value = 3784.6155
$ws.Range("K99").Value = 41673200
$ws.Range("L99").Value = 3784.6155
$ws.Range("M99").Value = -41671702
$ws.Range("N99").Value = -6780.6155
$ws.Range("H105").Value = 17606.584
$ws.Range("I105").Value = 38599.8
$ws.Range("J105").Value = 2611.4285
$ws.Range("K105").Value = 38599.8
$ws.Range("L105").Value = 2611.4285
$ws.Range("M105").Value = -36852.8
$ws.Range("N105").Value = -6105.4285
$ws.Range("H126").Value = 27032594
$ws.Range("I126").Value = 41673200
$ws.Range("J126").Value = 3784.6155
$ws.Range("K126").Value = 125019600
$ws.Range("L126").Value = 11353.8465
$ws.Range("M126").Value = -125017130
$ws.Range("N126").Value = -16293.8465
$ws.Range("H132").Value = 27780064
$ws.Range("I132").Value = 2288.8572
$ws.Range("J132").Value = 66668948
$ws.Range("K132").Value = 6866.571599999999
$ws.Range("L132").Value = 200006844
$ws.Range("M132").Value = -4336.571599999999
$ws.Range("N132").Value = -200011904

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H93").Value = 1005
$ws.Range("I93").Value = 1005
$ws.Range("K93").Value = 3015
$ws.Range("M93").Value = -1143
$ws.Range("H122").Value = 11578817
$ws.Range("I122").Value = 78125256
$ws.Range("J122").Value = 5523.478
$ws.Range("K122").Value = 703127304
$ws.Range("L122").Value = 49711.302
$ws.Range("M122").Value = -703124854
$ws.Range("N122").Value = -54611.302
$ws.Range("H131").Value = 854.08
$ws.Range("I131").Value = 450.75
$ws.Range("J131").Value = 909.0795000000001
$ws.Range("K131").Value = 1352.25
$ws.Range("L131").Value = 2727.2385
$ws.Range("M131").Value = 3687.75
$ws.Range("N131").Value = -12807.2385

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5924
$ws.Range("J70").Value = 6889.778
$ws.Range("L70").Value = 6889.778
$ws.Range("N70").Value = -7429.778
$ws.Range("H73").Value = 5924
$ws.Range("J73").Value = 6889.778
$ws.Range("L73").Value = 6889.778
$ws.Range("N73").Value = -8761.778
$ws.Range("H132").Value = 37802
$ws.Range("I132").Value = 33932.668
$ws.Range("J132").Value = 41671.332
$ws.Range("K132").Value = 101798.004
$ws.Range("L132").Value = 125013.996
$ws.Range("M132").Value = -99268.00399999999
$ws.Range("N132").Value = -130073.996

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 9683.076999999999
$ws.Range("I122").Value = 14237.5
$ws.Range("J122").Value = 2396
$ws.Range("K122").Value = 42712.5
$ws.Range("L122").Value = 7188
$ws.Range("M122").Value = -40262.5
$ws.Range("N122").Value = -12088
$ws.Range("H132").Value = 33339752
$ws.Range("I132").Value = 76926616
$ws.Range("J132").Value = 8621.177
$ws.Range("K132").Value = 230779848
$ws.Range("L132").Value = 25863.531
$ws.Range("M132").Value = -230777318
$ws.Range("N132").Value = -30923.531
$ws.Range("H136").Value = 36906644
$ws.Range("I136").Value = 14431839
$ws.Range("J136").Value = 142859300
$ws.Range("K136").Value = 43295517
$ws.Range("L136").Value = 428577900
$ws.Range("M136").Value = -43292967
$ws.Range("N136").Value = -428583000

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 36713.12
$ws.Range("I132").Value = 70500.266
$ws.Range("J132").Value = 8557.166999999999
$ws.Range("K132").Value = 211500.798
$ws.Range("L132").Value = 25671.501
$ws.Range("M132").Value = -208970.798
$ws.Range("N132").Value = -30731.501
$ws.Range("H136").Value = 8776165
$ws.Range("I136").Value = 12200124
$ws.Range("J136").Value = 2272.1875
$ws.Range("K136").Value = 36600372
$ws.Range("L136").Value = 6816.5625
$ws.Range("M136").Value = -36597822
$ws.Range("N136").Value = -11916.5625
